# Generate Report for Handback
# Updates the localization-status workbook: refreshes status text from
# "Ready for handoff" to "Handed back: in sync with en-US", refreshes the
# latest handback timestamps, clears the stale "version mismatch" error
# messages, and widens a couple of columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for both rows ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the status columns to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("K2").Value = "2016-11-09 00:25:08"
$zhcn.Range("K3").Value = "2016-11-09 00:25:08"

$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(16).ColumnWidth = 12.8

# --- de-de detail sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("K2").Value = "2016-11-09 00:25:25"
$dede.Range("K3").Value = "2016-11-09 00:25:25"

$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(16).ColumnWidth = 12.8
